$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo "Bewteen A and B." -> "Between A and B." for every row (column D, rows 1-70)
$ws.Range("D1:D70").Value = "Between A and B."

# Update the active selection to D8
$ws.Range("D8").Select()
